$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.517.37'
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').Value = '1.954.24'
$ws.Range('E3').Value = '  +0.68%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '244.06'
$ws.Range('D5').Style = "Normal"
$ws.Range('E6').Value = '  +0.37%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '58.09'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +1.40%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0790'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -6.32%  '
$ws.Range('E11').Value = '  -1.04%  '
$ws.Range('B12').Value = 'Polygon'
$ws.Range('C12').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.837'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +2.60%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '13.99'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +3.68%  '
$ws.Range('D14').Value = '2.241.48'
$ws.Range('E14').Value = '  +0.67%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '21.15'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.04%  '
$ws.Range('E16').Value = '  +1.75%  '
$ws.Range('D17').Value = '1.957.21'
$ws.Range('E17').Value = '  +1.14%  '
$ws.Range('D18').Value = '36.510.08'
$ws.Range('E18').Value = '  +0.50%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '69.68'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.40%  '
$ws.Range('E20').Value = '  -3.05%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '229.02'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.14%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.04'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.56%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.999'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.43'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +2.22%  '
$ws.Range('E25').Value = '  +2.52%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.140'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +6.39%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.13'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.66%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '160.01'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.58%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '19.33'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.23%  '
$ws.Range('E30').Value = '  +1.57%  '
$ws.Range('E31').Value = '  +4.61%  '
$ws.Range('E32').Value = '  +2.00%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0610'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -4.03%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.40'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +4.39%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.28'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +6.12%  '
$ws.Range('E36').Value = '  -0.03%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.44'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +12.71%  '
$ws.Range('E38').Value = '  -1.46%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.32'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -13.17%  '
$ws.Range('E40').Value = '  -0.31%  '
$ws.Range('E41').Value = '  +1.67%  '
$ws.Range('E42').Value = '  -0.18%  '
$ws.Range('E43').Value = '  -0.29%  '
$ws.Range('D44').Value = '1.371.09'
$ws.Range('E44').Value = '  +2.12%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '15.71'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.11%  '
$ws.Range('E46').Value = '  +0.61%  '
$ws.Range('E47').Value = '  -0.27%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '7.11'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.29%  '
$ws.Range('E49').Value = '  +0.15%  '
$ws.Range('D50').Value = '2.131.89'
$ws.Range('E50').Value = '  +0.68%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '43.91'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.59%  '
